$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("observed data")

# Copy row 88s formatting (date / time / datetime number formats) down
# across the new rows 89-109 so the freshly written cells pick up the same
# styles as the rest of the "observed data" table.
$ws.Range("A88:N88").Copy()
$ws.Range("A89:N109").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Row 89
$ws.Range("A89").Value = 44626
$ws.Range("B89").Value = "BriMac"
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 2
$ws.Range("E89").Value = "zone1"
$ws.Range("F89").Value = 0.51388888888888895
$ws.Range("G89").Value = 44626.513888888891
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("N89").Value = 21

# Row 90
$ws.Range("A90").Value = 44626
$ws.Range("B90").Value = "BriMac"
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 2
$ws.Range("E90").Value = "spur"
$ws.Range("F90").Value = 0.65972222222222221
$ws.Range("G90").Value = 44626.659722222219
$ws.Range("H90").Value = 2
$ws.Range("I90").Value = 1
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 1
$ws.Range("M90").Value = 0
$ws.Range("N90").Value = 11

# Row 91
$ws.Range("A91").Value = 44626
$ws.Range("B91").Value = "BriMac"
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 2
$ws.Range("E91").Value = "lot3"
$ws.Range("F91").Value = 0.63194444444444442
$ws.Range("G91").Value = 44626.631944444445
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 0
$ws.Range("N91").Value = 0

# Row 92
$ws.Range("A92").Value = 44619
$ws.Range("B92").Value = "BriMac"
$ws.Range("C92").Value = 10
$ws.Range("D92").Value = 4
$ws.Range("E92").Value = "zone1"
$ws.Range("F92").Value = 0.50694444444444442
$ws.Range("G92").Value = 44619.506944444445
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("N92").Value = 30

# Row 93
$ws.Range("A93").Value = 44619
$ws.Range("B93").Value = "BriMac"
$ws.Range("C93").Value = 10
$ws.Range("D93").Value = 4
$ws.Range("E93").Value = "spur"
$ws.Range("F93").Value = 0.54166666666666663
$ws.Range("G93").Value = 44619.541666666664
$ws.Range("H93").Value = 17
$ws.Range("I93").Value = 2
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("N93").Value = 17

# Row 94
$ws.Range("A94").Value = 44619
$ws.Range("B94").Value = "BriMac"
$ws.Range("C94").Value = 10
$ws.Range("D94").Value = 4
$ws.Range("E94").Value = "lot3"
$ws.Range("F94").Value = 0.5625
$ws.Range("G94").Value = 44619.5625
$ws.Range("H94").Value = 13
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 0
$ws.Range("N94").Value = 7

# Row 95
$ws.Range("A95").Value = 44612
$ws.Range("B95").Value = "BriMac"
$ws.Range("C95").Value = 50
$ws.Range("D95").Value = 2
$ws.Range("E95").Value = "zone1"
$ws.Range("F95").Value = 0.51388888888888895
$ws.Range("G95").Value = 44612.513888888891
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 1
$ws.Range("M95").Value = 0
$ws.Range("N95").Value = 30
$ws.Range("O95").Value = "people too numerous to count; est"

# Row 96
$ws.Range("A96").Value = 44612
$ws.Range("B96").Value = "BriMac"
$ws.Range("C96").Value = 50
$ws.Range("D96").Value = 2
$ws.Range("E96").Value = "spur"
$ws.Range("F96").Value = 0.625
$ws.Range("G96").Value = 44612.625
$ws.Range("H96").Value = 5
$ws.Range("I96").Value = 2
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1
$ws.Range("M96").Value = 0
$ws.Range("N96").Value = 11

# Row 97
$ws.Range("A97").Value = 44612
$ws.Range("B97").Value = "BriMac"
$ws.Range("C97").Value = 50
$ws.Range("D97").Value = 2
$ws.Range("E97").Value = "lot3"
$ws.Range("F97").Value = 0.64583333333333337
$ws.Range("G97").Value = 44612.645833333336
$ws.Range("H97").Value = 4
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("N97").Value = 8

# Row 98
$ws.Range("A98").Value = 44584
$ws.Range("B98").Value = "BriMac"
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 2
$ws.Range("E98").Value = "zone1"
$ws.Range("F98").Value = 0.53125
$ws.Range("G98").Value = 44584.53125
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("N98").Value = 15

# Row 99
$ws.Range("A99").Value = 44584
$ws.Range("B99").Value = "BriMac"
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 2
$ws.Range("E99").Value = "spur"
$ws.Range("F99").Value = 0.63194444444444442
$ws.Range("G99").Value = 44584.631944444445
$ws.Range("H99").Value = 13
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("N99").Value = 11

# Row 100
$ws.Range("A100").Value = 44584
$ws.Range("B100").Value = "BriMac"
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 2
$ws.Range("E100").Value = "lot3"
$ws.Range("F100").Value = 0.65625
$ws.Range("G100").Value = 44584.65625
$ws.Range("H100").Value = 8
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 0
$ws.Range("N100").Value = 5

# Row 101
$ws.Range("A101").Value = 44591
$ws.Range("B101").Value = "BriMac"
$ws.Range("C101").Value = 20
$ws.Range("D101").Value = 3
$ws.Range("E101").Value = "zone1"
$ws.Range("F101").Value = 0.51041666666666663
$ws.Range("G101").Value = 44591.510416666664
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("O101").Value = "people too numerous to count"

# Row 102
$ws.Range("A102").Value = 44591
$ws.Range("B102").Value = "BriMac"
$ws.Range("C102").Value = 20
$ws.Range("D102").Value = 3
$ws.Range("E102").Value = "spur"
$ws.Range("F102").Value = 0.55208333333333337
$ws.Range("G102").Value = 44591.552083333336
$ws.Range("H102").Value = 20
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 0
$ws.Range("N102").Value = 18

# Row 103
$ws.Range("A103").Value = 44591
$ws.Range("B103").Value = "BriMac"
$ws.Range("C103").Value = 20
$ws.Range("D103").Value = 3
$ws.Range("E103").Value = "lot3"
$ws.Range("F103").Value = 0.57291666666666663
$ws.Range("G103").Value = 44591.572916666664
$ws.Range("H103").Value = 14
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 0
$ws.Range("N103").Value = 0

# Row 104
$ws.Range("A104").Value = 44597
$ws.Range("B104").Value = "KatGer"
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 2
$ws.Range("E104").Value = "zone1"
$ws.Range("F104").Value = 0.56319444444444444
$ws.Range("G104").Value = 44597.563194444447
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 0
$ws.Range("N104").Value = 4

# Row 105
$ws.Range("A105").Value = 44597
$ws.Range("B105").Value = "KatGer"
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 2
$ws.Range("E105").Value = "spur"
$ws.Range("F105").Value = 0.5805555555555556
$ws.Range("G105").Value = 44597.580555555556
$ws.Range("H105").Value = 15
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 1
$ws.Range("M105").Value = 0
$ws.Range("N105").Value = 0

# Row 106
$ws.Range("A106").Value = 44597
$ws.Range("B106").Value = "KatGer"
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 2
$ws.Range("E106").Value = "lot3"
$ws.Range("F106").Value = 0.60277777777777775
$ws.Range("G106").Value = 44597.602777777778
$ws.Range("H106").Value = 16
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = 0
$ws.Range("N106").Value = 5

# Row 107
$ws.Range("A107").Value = 44604
$ws.Range("B107").Value = "KatGer"
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 1
$ws.Range("E107").Value = "zone1"
$ws.Range("F107").Value = 0.55208333333333337
$ws.Range("G107").Value = 44604.552083333336
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 0
$ws.Range("N107").Value = 11

# Row 108
$ws.Range("A108").Value = 44604
$ws.Range("B108").Value = "KatGer"
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 1
$ws.Range("E108").Value = "spur"
$ws.Range("F108").Value = 0.57986111111111105
$ws.Range("G108").Value = 44604.579861111109
$ws.Range("H108").Value = 16
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 1
$ws.Range("M108").Value = 0
$ws.Range("N108").Value = 4

# Row 109
$ws.Range("A109").Value = 44604
$ws.Range("B109").Value = "KatGer"
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 1
$ws.Range("E109").Value = "lot3"
$ws.Range("F109").Value = 0.61458333333333337
$ws.Range("G109").Value = 44604.614583333336
$ws.Range("H109").Value = 10
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = 0
$ws.Range("N109").Value = 2
